$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns (values like "1.00" or
# "97.144.78" must stay literal text, not be coerced to numbers/dates)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "97.144.78"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.711.82"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "236.61"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "657.15"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.432"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").Value = "3.708.42"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "0.0000318"
$ws.Range("E12").Value = "  +17.56%  "
$ws.Range("D13").Value = "44.66"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "4.406.04"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "96.816.04"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "8.97"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "3.717.96"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "13.01"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "18.84"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "0.504"
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("D23").Value = "523.46"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "3.46"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "0.0000224"
$ws.Range("E25").Value = "  +9.76%  "
$ws.Range("D26").Value = "6.88"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").Value = "106.19"
$ws.Range("E27").Value = "  +3.68%  "
$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").Value = "0.191"
$ws.Range("E28").Value = "  +13.80%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.933.17"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "13.55"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "12.71"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").Value = "3.03"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  -3.80%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "32.41"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").Value = "637.49"
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("D39").Value = "0.590"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").Value = "8.74"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D42").Value = "0.165"
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "40.62"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "6.71"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.486"
$ws.Range("E45").Value = "  +13.32%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "2.01"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.966"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "0.0455"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "2.36"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").Value = "23.62"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "8.63"
$ws.Range("E51").Value = "  -0.80%  "
